$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style incl. borders) from row 45 down into the new rows 46:51
$ws.Range("A45:B45").Copy()
$ws.Range("A46:B51").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New contributor rows
$ws.Range("A46").Value = "Christina Wang"
$ws.Range("B46").Value = "NYU"
$ws.Range("A47").Value = "Tim Levengood, MPH"
$ws.Range("B47").Value = "Boston University"
$ws.Range("A48").Value = "María Jordán P."
$ws.Range("B48").Value = "Boston University"
$ws.Range("A49").Value = "Joel Earlywine"
$ws.Range("B49").Value = "Boston University"
$ws.Range("A50").Value = "Grace Yoon"
$ws.Range("B50").Value = "Boston University"
$ws.Range("A51").Value = "Emma Chistolini"
$ws.Range("B51").Value = "Colby College"

# Uniform row height (15) for the contributor rows, matching the author's re-format
$ws.Range("A15:B51").EntireRow.RowHeight = 15

# Selection / scroll position left by the editing session
$ws.Range("A53").Select()
